$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The underlying word list was re-generated by the NLP text-processing step and
# came back in a different order. Column A is updated cell-by-cell so every row
# shows the regenerated word for that position; column B (the frequency counts)
# is left untouched, exactly as it was before the edit.
$ws.Range("A17").Value = "колеса"
$ws.Range("A18").Value = "Крымскую соль"
$ws.Range("A19").Value = "парча"
$ws.Range("A20").Value = "говядина"
$ws.Range("A21").Value = "сено"
$ws.Range("A24").Value = "чулок"
$ws.Range("A25").Value = "выбойка"
$ws.Range("A26").Value = "сахар"
$ws.Range("A27").Value = "шелк"
$ws.Range("A31").Value = "сани"
$ws.Range("A34").Value = "коса"
$ws.Range("A35").Value = "платок"
$ws.Range("A36").Value = "рогожа"
$ws.Range("A37").Value = "замок"
$ws.Range("A38").Value = "обод"
$ws.Range("A39").Value = "овца"
$ws.Range("A41").Value = "конь"
$ws.Range("A42").Value = "веревка"
$ws.Range("A44").Value = "ром"
$ws.Range("A45").Value = "нитка"
$ws.Range("A46").Value = "скотский кожа"
$ws.Range("A47").Value = "гумми"
$ws.Range("A48").Value = "брусья"
$ws.Range("A49").Value = "котел"
$ws.Range("A50").Value = "покроми"
$ws.Range("A51").Value = "сковорода"
$ws.Range("A52").Value = "дуга"
$ws.Range("A53").Value = "хомут"
$ws.Range("A54").Value = "бечева"
$ws.Range("A55").Value = "сосуд"
$ws.Range("A56").Value = "роза"
